$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rate entry appended as row 17
$cellA = $ws.Cells.Item(17, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-09-11"
$cellA.Style = "Normal"

$ws.Range("B17").Value = "15:20:00"
$ws.Range("C17").Value = "1.00 EUR = 1667.5922 ARS"
